$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the A column SPI pins (bottom-up so shared-string order matches
# the order the labels were typed in: CS, MOSI, MISO, CLK)
$ws.Range("A8").Value = "SPI CS"
$ws.Range("A7").Value = "SPI MOSI"
$ws.Range("A6").Value = "SPI MISO"
$ws.Range("A5").Value = "SPI CLK"

# Add the new Lift Pul / Lift Dir labels in the N column, matching the
# existing left-aligned / highlighted style used by N6:N7
$ws.Range("N4").Value = "Lift Pul"
$ws.Range("N4").HorizontalAlignment = -4131
$ws.Range("N5").Value = "Lift Dir"
$ws.Range("N5").HorizontalAlignment = -4131

# Move the selection to O9
$ws.Range("O9").Select()
